$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 4631261.5
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 5557314
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 5557314
$ws.Range("M43").Value = -931
$ws.Range("N43").Value = -5557452
$ws.Range("H58").Value = 859.1579
$ws.Range("I58").Value = 440.83334
$ws.Range("J58").Value = 1576.2858
$ws.Range("K58").Value = 1322.50002
$ws.Range("L58").Value = 4728.857400000001
$ws.Range("M58").Value = -1172.50002
$ws.Range("N58").Value = -5028.857400000001
$ws.Range("H138").Value = 1951.63
$ws.Range("J138").Value = 2206.0264
$ws.Range("L138").Value = 6618.0792
$ws.Range("N138").Value = -16898.0792
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1508.2307
$ws.Range("I2").Value = 1067.6666
$ws.Range("K2").Value = 1067.6666
$ws.Range("M2").Value = -954.6666
$ws.Range("H45").Value = 1046.6666
$ws.Range("I45").Value = 920.6667
$ws.Range("K45").Value = 920.6667
$ws.Range("M45").Value = -543.6667
$ws.Range("H61").Value = 27778910
$ws.Range("I61").Value = 30303972
$ws.Range("K61").Value = 30303972
$ws.Range("M61").Value = -30303760
$ws.Range("H74").Value = 1996
$ws.Range("I74").Value = 1368
$ws.Range("J74").Value = 3880
$ws.Range("K74").Value = 1368
$ws.Range("L74").Value = 3880
$ws.Range("M74").Value = -494
$ws.Range("N74").Value = -5628
$ws.Range("H77").Value = 1996
$ws.Range("I77").Value = 1368
$ws.Range("J77").Value = 3880
$ws.Range("K77").Value = 6840
$ws.Range("L77").Value = 19400
$ws.Range("M77").Value = -2472
$ws.Range("N77").Value = -28136
$ws.Range("H110").Value = 850.3
$ws.Range("I110").Value = 387.77777
$ws.Range("K110").Value = 387.77777
$ws.Range("M110").Value = 1657.22223
$ws.Range("H116").Value = 1508.2307
$ws.Range("I116").Value = 1067.6666
$ws.Range("K116").Value = 1067.6666
$ws.Range("M116").Value = 1226.3334
$ws.Range("H122").Value = 3288
$ws.Range("I122").Value = 3334.8333
$ws.Range("K122").Value = 10004.4999
$ws.Range("M122").Value = -7554.499899999999
$ws.Range("H132").Value = 3408.9697
$ws.Range("I132").Value = 2747.6667
$ws.Range("K132").Value = 8243.000100000001
$ws.Range("M132").Value = -5713.000100000001
$ws.Range("H135").Value = 80429
$ws.Range("J135").Value = 80429
$ws.Range("L135").Value = 80429
$ws.Range("N135").Value = -90569
$ws.Range("H136").Value = 27778910
$ws.Range("I136").Value = 30303972
$ws.Range("K136").Value = 90911916
$ws.Range("M136").Value = -90909366
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1508.2307
$ws.Range("I3").Value = 1067.6666
$ws.Range("K3").Value = 1067.6666
$ws.Range("M3").Value = -953.6666
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H81").Value = 14977
$ws.Range("J81").Value = 14977
$ws.Range("L81").Value = 14977
$ws.Range("N81").Value = -17099
$ws.Range("H84").Value = 14977
$ws.Range("J84").Value = 14977
$ws.Range("L84").Value = 44931
$ws.Range("N84").Value = -55539
$ws.Range("H107").Value = 917.34485
$ws.Range("I107").Value = 695.4167
$ws.Range("K107").Value = 695.4167
$ws.Range("M107").Value = 1224.5833
$ws.Range("H112").Value = 134449.5
$ws.Range("J112").Value = 134449.5
$ws.Range("L112").Value = 134449.5
$ws.Range("N112").Value = -137403.5
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H116").Value = 69999.664
$ws.Range("J116").Value = 69999.664
$ws.Range("L116").Value = 69999.664
$ws.Range("N116").Value = -79177.664
$ws.Range("H119").Value = 29880.5
$ws.Range("J119").Value = 29880.5
$ws.Range("L119").Value = 29880.5
$ws.Range("N119").Value = -39556.5
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("M121").ClearContents()
$ws.Range("H134").Value = 1704.7
$ws.Range("I134").Value = 1430.875
$ws.Range("J134").Value = 2800
$ws.Range("K134").Value = 4292.625
$ws.Range("L134").Value = 8400
$ws.Range("M134").Value = -1757.625
$ws.Range("N134").Value = -13470
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 1750
$ws.Range("I45").Value = 1750
$ws.Range("K45").Value = 1750
$ws.Range("M45").Value = -1157
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("H122").Value = 1007.625
$ws.Range("I122").Value = 1030.6666
$ws.Range("J122").Value = 978
$ws.Range("K122").Value = 3091.9998
$ws.Range("L122").Value = 2934
$ws.Range("M122").Value = -641.9998000000001
$ws.Range("N122").Value = -7834
$ws.Range("H134").Value = 17858754
$ws.Range("I134").Value = 1626.4286
$ws.Range("J134").Value = 71430140
$ws.Range("K134").Value = 4879.2858
$ws.Range("L134").Value = 214290420
$ws.Range("M134").Value = -2344.2858
$ws.Range("N134").Value = -214295490
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 130.375
$ws.Range("I12").Value = 262.5
$ws.Range("J12").Value = 86.333336
$ws.Range("K12").Value = 787.5
$ws.Range("L12").Value = 259.000008
$ws.Range("M12").Value = -614.5
$ws.Range("N12").Value = -605.000008
$ws.Range("H57").Value = 750
$ws.Range("I57").Value = 750
$ws.Range("K57").Value = 2250
$ws.Range("M57").Value = -1691
$ws.Range("H131").Value = 20439558
$ws.Range("J131").Value = 39312.28
$ws.Range("L131").Value = 117936.84
$ws.Range("N131").Value = -128016.84
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1041.4
$ws.Range("J97").Value = 1002.3333
$ws.Range("L97").Value = 1002.3333
$ws.Range("N97").Value = -1994.3333
$ws.Range("H102").Value = 1599.1852
$ws.Range("I102").Value = 1543.7222
$ws.Range("J102").Value = 1710.1111
$ws.Range("K102").Value = 1543.7222
$ws.Range("L102").Value = 1710.1111
$ws.Range("M102").Value = 78.27780000000007
$ws.Range("N102").Value = -4954.1111
$ws.Range("H122").Value = 3750
$ws.Range("I122").Value = 3750
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 11250
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8800
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 3741.5217
$ws.Range("I132").Value = 3826.3845
$ws.Range("J132").Value = 3631.2
$ws.Range("K132").Value = 11479.1535
$ws.Range("L132").Value = 10893.6
$ws.Range("M132").Value = -8949.1535
$ws.Range("N132").Value = -15953.6
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3731.1765
$ws.Range("I46").Value = 916.6667
$ws.Range("J46").Value = 5266.364
$ws.Range("K46").Value = 916.6667
$ws.Range("L46").Value = 5266.364
$ws.Range("M46").Value = -728.6667
$ws.Range("N46").Value = -5642.364
$ws.Range("H133").Value = 49613
$ws.Range("J133").Value = 49613
$ws.Range("L133").Value = 49613
$ws.Range("N133").Value = -54673
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 50000
$ws.Range("J16").Value = 50000
$ws.Range("L16").Value = 50000
$ws.Range("N16").Value = -50584
$ws.Range("H107").Value = 603.9032
$ws.Range("I107").Value = 380.53333
$ws.Range("J107").Value = 813.3125
$ws.Range("K107").Value = 1141.59999
$ws.Range("L107").Value = 2439.9375
$ws.Range("M107").Value = 778.4000100000001
$ws.Range("N107").Value = -6279.9375
$ws.Range("H132").Value = 2003.0869
$ws.Range("I132").Value = 1826.7354
$ws.Range("J132").Value = 2502.75
$ws.Range("K132").Value = 5480.206200000001
$ws.Range("L132").Value = 7508.25
$ws.Range("M132").Value = -2950.206200000001
$ws.Range("N132").Value = -12568.25
